$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header row (row 1) ---
# C1 used to read "Boolean"; it's now split into two explicit boolean
# columns plus a new "Input" column.
$ws.Range("C1").Value = "Boolean 1"
$ws.Range("D1").Value = "Bollean 0"
$ws.Range("E1").Value = "Input"

# --- Update existing data row (row 2) ---
$ws.Range("B2").Value = "Please enter the invoice split count less than 7"

# --- Add new data row 3 (TC_002) ---
$ws.Range("A3").Value = "TC_002"
$ws.Range("C3").Value = $true

# --- Add new data row 4 (TC_003 TC_005 TC_006) ---
$ws.Range("A4").Value = "TC_003 TC_005 TC_006"
$ws.Range("C4").Value = $true
$ws.Range("D4").Value = $false
$ws.Range("E4").Value = 123456789

# --- Column sizing to fit the new content ---
$ws.Columns.Item(1).ColumnWidth = 21.8333333333333
$ws.Columns.Item(5).ColumnWidth = 9.2

# --- Selection: whole of column A is selected, as if its header was clicked ---
$ws.Range("A1:A1048576").Select() | Out-Null
